# Auto-generated Excel COM-interop script to apply market-data refresh edits
# described by the commit "chore: update Sheets via scheduled runner".
# For each affected sheet/row, the currentAveragePrice* / Leve*Price* / LeveProfit*
# columns (H:N) are updated to the newly recomputed values. Cells that are added
# or removed by the update are handled explicitly (SetValue / ClearContents).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2835.625
$ws.Range("J62").Value = 3984.5
$ws.Range("L62").Value = 3984.5
$ws.Range("N62").Value = -5232.5
$ws.Range("H65").Value = 2835.625
$ws.Range("J65").Value = 3984.5
$ws.Range("L65").Value = 19922.5
$ws.Range("N65").Value = -26162.5
$ws.Range("H98").Value = 1698
$ws.Range("I98").Value = 1410.25
$ws.Range("K98").Value = 1410.25
$ws.Range("M98").Value = 87.75
$ws.Range("H112").Value = 7016.923
$ws.Range("I112").Value = 70100
$ws.Range("J112").Value = 1760
$ws.Range("K112").Value = 210300
$ws.Range("L112").Value = 5280
$ws.Range("M112").Value = -209192
$ws.Range("N112").Value = -7496
$ws.Range("H122").Value = 1698
$ws.Range("I122").Value = 1410.25
$ws.Range("K122").Value = 4230.75
$ws.Range("M122").Value = -1780.75
$ws.Range("H129").Value = 931.5208
$ws.Range("I129").Value = 315.83334
$ws.Range("J129").Value = 1019.4762
$ws.Range("K129").Value = 947.5000200000001
$ws.Range("L129").Value = 3058.4286
$ws.Range("M129").Value = 4052.49998
$ws.Range("N129").Value = -13058.4286
$ws.Range("H137").Value = 669542.7
$ws.Range("I137").Value = 3300.1052
$ws.Range("J137").Value = 1077884.9
$ws.Range("K137").Value = 9900.3156
$ws.Range("L137").Value = 3233654.7
$ws.Range("M137").Value = -7350.3156
$ws.Range("N137").Value = -3238754.7
$ws.Range("H138").Value = 4728
$ws.Range("I138").Value = 2080.4375
$ws.Range("J138").Value = 6011.6665
$ws.Range("K138").Value = 6241.3125
$ws.Range("L138").Value = 18034.9995
$ws.Range("M138").Value = -1101.3125
$ws.Range("N138").Value = -28314.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27530.666
$ws.Range("I32").Value = 28878.46
$ws.Range("J32").Value = 10009.333
$ws.Range("K32").Value = 28878.46
$ws.Range("L32").Value = 10009.333
$ws.Range("M32").Value = -28591.46
$ws.Range("N32").Value = -10583.333
$ws.Range("H61").Value = 5290.095
$ws.Range("I61").Value = 3966.3547
$ws.Range("J61").Value = 9020.637000000001
$ws.Range("K61").Value = 3966.3547
$ws.Range("L61").Value = 9020.637000000001
$ws.Range("M61").Value = -3754.3547
$ws.Range("N61").Value = -9444.637000000001
$ws.Range("H74").Value = 4133.4326
$ws.Range("I74").Value = 1453.2069
$ws.Range("J74").Value = 13849.25
$ws.Range("K74").Value = 1453.2069
$ws.Range("L74").Value = 13849.25
$ws.Range("M74").Value = -579.2068999999999
$ws.Range("N74").Value = -15597.25
$ws.Range("H77").Value = 4133.4326
$ws.Range("I77").Value = 1453.2069
$ws.Range("J77").Value = 13849.25
$ws.Range("K77").Value = 7266.0345
$ws.Range("L77").Value = 69246.25
$ws.Range("M77").Value = -2898.0345
$ws.Range("N77").Value = -77982.25
$ws.Range("H101").Value = 45602
$ws.Range("J101").Value = 45602
$ws.Range("L101").Value = 45602
$ws.Range("N101").Value = -52092
$ws.Range("H102").Value = 4116.5
$ws.Range("I102").Value = 3974.75
$ws.Range("K102").Value = 3974.75
$ws.Range("M102").Value = -2352.75
$ws.Range("H132").Value = 1992.4615
$ws.Range("I132").Value = 1806.4897
$ws.Range("J132").Value = 2562
$ws.Range("K132").Value = 5419.4691
$ws.Range("L132").Value = 7686
$ws.Range("M132").Value = -2889.4691
$ws.Range("N132").Value = -12746
$ws.Range("H133").Value = 43168.555
$ws.Range("J133").Value = 43168.555
$ws.Range("L133").Value = 43168.555
$ws.Range("N133").Value = -48228.555
$ws.Range("H136").Value = 5290.095
$ws.Range("I136").Value = 3966.3547
$ws.Range("J136").Value = 9020.637000000001
$ws.Range("K136").Value = 11899.0641
$ws.Range("L136").Value = 27061.911
$ws.Range("M136").Value = -9349.0641
$ws.Range("N136").Value = -32161.911

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11600
$ws.Range("I105").Value = 14666.667
$ws.Range("J105").Value = 7000
$ws.Range("K105").Value = 14666.667
$ws.Range("L105").Value = 7000
$ws.Range("M105").Value = -12919.667
$ws.Range("N105").Value = -10494
$ws.Range("H134").Value = 2759.5334
$ws.Range("I134").Value = 3030.2307
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 9090.6921
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -6555.6921
$ws.Range("N134").Value = -8070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 10000
$ws.Range("I32").Value = 10000
$ws.Range("K32").Value = 10000
$ws.Range("M32").Value = -9684
$ws.Range("H58").Value = 1980195.9
$ws.Range("I58").Value = 2675104.8
$ws.Range("J58").Value = 11287.833
$ws.Range("K58").Value = 2675104.8
$ws.Range("L58").Value = 11287.833
$ws.Range("M58").Value = -2674901.8
$ws.Range("N58").Value = -11693.833
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()
$ws.Range("H106").Value = 42900
$ws.Range("J106").Value = 42900
$ws.Range("L106").Value = 42900
$ws.Range("N106").Value = -45424
$ws.Range("H132").Value = 2251.8281
$ws.Range("I132").Value = 2038.2667
$ws.Range("J132").Value = 2757.6316
$ws.Range("K132").Value = 6114.800099999999
$ws.Range("L132").Value = 8272.8948
$ws.Range("M132").Value = -3584.800099999999
$ws.Range("N132").Value = -13332.8948
$ws.Range("H136").Value = 1980195.9
$ws.Range("I136").Value = 2675104.8
$ws.Range("J136").Value = 11287.833
$ws.Range("K136").Value = 8025314.399999999
$ws.Range("L136").Value = 33863.499
$ws.Range("M136").Value = -8022764.399999999
$ws.Range("N136").Value = -38963.499

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2733.9404
$ws.Range("J68").Value = 4486.1665
$ws.Range("L68").Value = 13458.4995
$ws.Range("N68").Value = -15080.4995
$ws.Range("H71").Value = 2733.9404
$ws.Range("J71").Value = 4486.1665
$ws.Range("L71").Value = 40375.4985
$ws.Range("N71").Value = -48487.4985
$ws.Range("H107").Value = 1075.553
$ws.Range("I107").Value = 588.8333
$ws.Range("J107").Value = 1341.0364
$ws.Range("K107").Value = 1766.4999
$ws.Range("L107").Value = 4023.1092
$ws.Range("M107").Value = 153.5001
$ws.Range("N107").Value = -7863.1092
$ws.Range("H122").Value = 1149.7142
$ws.Range("I122").Value = 734
$ws.Range("J122").Value = 1279.625
$ws.Range("K122").Value = 6606
$ws.Range("L122").Value = 11516.625
$ws.Range("M122").Value = -4156
$ws.Range("N122").Value = -16416.625
$ws.Range("H131").Value = 44982
$ws.Range("I131").Value = 1719
$ws.Range("J131").Value = 88245
$ws.Range("K131").Value = 5157
$ws.Range("L131").Value = 264735
$ws.Range("M131").Value = -117
$ws.Range("N131").Value = -274815

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1990.909
$ws.Range("I113").Value = 1990.909
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1990.909
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 179.0909999999999
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 33871.312
$ws.Range("I132").Value = 47234.637
$ws.Range("J132").Value = 4472
$ws.Range("K132").Value = 141703.911
$ws.Range("L132").Value = 13416
$ws.Range("M132").Value = -139173.911
$ws.Range("N132").Value = -18476

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 37965.547
$ws.Range("I61").Value = 36093.24
$ws.Range("K61").Value = 36093.24
$ws.Range("M61").Value = -35891.24
$ws.Range("H95").Value = 4344
$ws.Range("J95").Value = 4344
$ws.Range("L95").Value = 4344
$ws.Range("N95").Value = -9836
$ws.Range("H100").Value = 5571.4287
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459
$ws.Range("H113").Value = 37965.547
$ws.Range("I113").Value = 36093.24
$ws.Range("K113").Value = 36093.24
$ws.Range("M113").Value = -33923.24
$ws.Range("H122").Value = 6543
$ws.Range("I122").Value = 6313.3335
$ws.Range("J122").Value = 7921
$ws.Range("K122").Value = 18940.0005
$ws.Range("L122").Value = 23763
$ws.Range("M122").Value = -16490.0005
$ws.Range("N122").Value = -28663
$ws.Range("H132").Value = 9841.409
$ws.Range("I132").Value = 13579.071
$ws.Range("J132").Value = 3300.5
$ws.Range("K132").Value = 40737.213
$ws.Range("L132").Value = 9901.5
$ws.Range("M132").Value = -38207.213
$ws.Range("N132").Value = -14961.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 52200
$ws.Range("I2").Value = 126000
$ws.Range("K2").Value = 126000
$ws.Range("M2").Value = -125888
$ws.Range("H81").Value = 40003500
$ws.Range("I81").Value = 3749
$ws.Range("J81").Value = 66670000
$ws.Range("K81").Value = 7498
$ws.Range("L81").Value = 133340000
$ws.Range("M81").Value = -6437
$ws.Range("N81").Value = -133342122
$ws.Range("H84").Value = 40003500
$ws.Range("I84").Value = 3749
$ws.Range("J84").Value = 66670000
$ws.Range("K84").Value = 37490
$ws.Range("L84").Value = 666700000
$ws.Range("M84").Value = -32186
$ws.Range("N84").Value = -666710608
$ws.Range("H113").Value = 751.4400000000001
$ws.Range("I113").Value = 473.8889
$ws.Range("J113").Value = 907.5625
$ws.Range("K113").Value = 1421.6667
$ws.Range("L113").Value = 2722.6875
$ws.Range("M113").Value = 748.3333
$ws.Range("N113").Value = -7062.6875
$ws.Range("H132").Value = 2587.9678
$ws.Range("I132").Value = 2491.9565
$ws.Range("J132").Value = 2864
$ws.Range("K132").Value = 7475.869499999999
$ws.Range("L132").Value = 8592
$ws.Range("M132").Value = -4945.869499999999
$ws.Range("N132").Value = -13652
